$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-28 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-29 Friday", 2) | Out-Null
$d.Content.Find.Execute("599×3=1797", $true, $false, $false, $false, $false, $true, 1, $false, "948×9=8532", 2) | Out-Null
$d.Content.Find.Execute("613×4=2452", $true, $false, $false, $false, $false, $true, 1, $false, "433×9=3897", 2) | Out-Null
$d.Content.Find.Execute("307×3=921", $true, $false, $false, $false, $false, $true, 1, $false, "171×2=342", 2) | Out-Null
$d.Content.Find.Execute("771×8=6168", $true, $false, $false, $false, $false, $true, 1, $false, "419×5=2095", 2) | Out-Null
$d.Content.Find.Execute("762×4=3048", $true, $false, $false, $false, $false, $true, 1, $false, "927×9=8343", 2) | Out-Null
$d.Content.Find.Execute("317×9=2853", $true, $false, $false, $false, $false, $true, 1, $false, "601×3=1803", 2) | Out-Null
$d.Content.Find.Execute("119×7=833", $true, $false, $false, $false, $false, $true, 1, $false, "389×7=2723", 2) | Out-Null
$d.Content.Find.Execute("537×7=3759", $true, $false, $false, $false, $false, $true, 1, $false, "921×7=6447", 2) | Out-Null
$d.Content.Find.Execute("887×2=1774", $true, $false, $false, $false, $false, $true, 1, $false, "463×8=3704", 2) | Out-Null
$d.Content.Find.Execute("453×3=1359", $true, $false, $false, $false, $false, $true, 1, $false, "119×4=476", 2) | Out-Null
$d.Content.Find.Execute("638×7=4466", $true, $false, $false, $false, $false, $true, 1, $false, "957×4=3828", 2) | Out-Null
$d.Content.Find.Execute("749×5=3745", $true, $false, $false, $false, $false, $true, 1, $false, "375×8=3000", 2) | Out-Null
$d.Content.Find.Execute("263×7=1841", $true, $false, $false, $false, $false, $true, 1, $false, "248×8=1984", 2) | Out-Null
$d.Content.Find.Execute("255×2=510", $true, $false, $false, $false, $false, $true, 1, $false, "536×2=1072", 2) | Out-Null
$d.Content.Find.Execute("632×9=5688", $true, $false, $false, $false, $false, $true, 1, $false, "351×3=1053", 2) | Out-Null
$d.Content.Find.Execute("295×3=885", $true, $false, $false, $false, $false, $true, 1, $false, "208×3=624", 2) | Out-Null
$d.Content.Find.Execute("839×8=6712", $true, $false, $false, $false, $false, $true, 1, $false, "130×4=520", 2) | Out-Null
$d.Content.Find.Execute("859×5=4295", $true, $false, $false, $false, $false, $true, 1, $false, "708×6=4248", 2) | Out-Null
$d.Content.Find.Execute("483×9=4347", $true, $false, $false, $false, $false, $true, 1, $false, "258×5=1290", 2) | Out-Null
$d.Content.Find.Execute("723×9=6507", $true, $false, $false, $false, $false, $true, 1, $false, "971×5=4855", 2) | Out-Null
$d.Content.Find.Execute("850×3=2550", $true, $false, $false, $false, $false, $true, 1, $false, "710×3=2130", 2) | Out-Null
$d.Content.Find.Execute("803×3=2409", $true, $false, $false, $false, $false, $true, 1, $false, "868×5=4340", 2) | Out-Null
$d.Content.Find.Execute("876×3=2628", $true, $false, $false, $false, $false, $true, 1, $false, "331×6=1986", 2) | Out-Null
$d.Content.Find.Execute("954×4=3816", $true, $false, $false, $false, $false, $true, 1, $false, "406×2=812", 2) | Out-Null
$d.Content.Find.Execute("578×2=1156", $true, $false, $false, $false, $false, $true, 1, $false, "835×5=4175", 2) | Out-Null
